$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 23:38"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 8726248
$ws.Range("C4").Value = 61067
$ws.Range("D4").Value = 5690069
$ws.Range("E4").Value = 2807034
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 764
$ws.Range("H4").Value = 229145

# Row 6: Brasil
$ws.Range("A6").Value = "Brasil"
$ws.Range("B6").Value = 5352935
$ws.Range("C6").Value = 20301
$ws.Range("D6").Value = 4785297
$ws.Range("E6").Value = 411169
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 507
$ws.Range("H6").Value = 156469

# Row 129: Trinidad yTobago
$ws.Range("A129").Value = "Trinidad yTobago"
$ws.Range("B129").Value = 5487
$ws.Range("C129").Value = 41
$ws.Range("D129").Value = 3945
$ws.Range("E129").Value = 1438
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 104

# Row 139: Ruanda
$ws.Range("A139").Value = "Ruanda"
$ws.Range("B139").Value = 5052
$ws.Range("C139").Value = 35
$ws.Range("D139").Value = 4806
$ws.Range("E139").Value = 212
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 34

# Row 152: Republica de Chipre
$ws.Range("A152").Value = "Republica de Chipre"
$ws.Range("B152").Value = 3314
$ws.Range("C152").Value = 160
$ws.Range("D152").Value = 1444
$ws.Range("E152").Value = 1845
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 25

# Row 162: Yemen
$ws.Range("A162").Value = "Yemen"
$ws.Range("B162").Value = 2060
$ws.Range("C162").Value = 3
$ws.Range("D162").Value = 1354
$ws.Range("E162").Value = 107
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 2
$ws.Range("H162").Value = 599

# Row 199: San Vicente y las Granadinas
$ws.Range("A199").Value = "San Vicente y las Granadinas"
$ws.Range("B199").Value = 73
$ws.Range("C199").Value = 5
$ws.Range("D199").Value = 64
$ws.Range("E199").Value = 9
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

# Row 200: Islas Virgenes Britanicas
$ws.Range("A200").Value = "Islas Virgenes Britanicas"
$ws.Range("B200").Value = 71
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 70
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 1

# Row 216: Islas Malvinas
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 13
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0

# Row 217: Montserrat
$ws.Range("A217").Value = "Montserrat"
$ws.Range("B217").Value = 13
$ws.Range("C217").Value = 0
$ws.Range("D217").Value = 12
$ws.Range("E217").Value = 0
$ws.Range("F217").Value = 0
$ws.Range("G217").Value = 0
$ws.Range("H217").Value = 1
